$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 6.755097
$ws.Range("N2").Value = 20.265291
$ws.Range("O2").Value = 0.262181130417721
$ws.Range("P2").Value = 0.262181130417721
$ws.Range("Q2").Value = 7.134839281253
$ws.Range("R2").Value = 64.213553531277
$ws.Range("S2").Value = 0.262181130417721
$ws.Range("T2").Value = 0.262181130417721

# Row 3
$ws.Range("O3").Value = 0.179964029239562
$ws.Range("P3").Value = 0.179964029239562
$ws.Range("S3").Value = 0.179964029239562
$ws.Range("T3").Value = 0.179964029239562

# Row 4
$ws.Range("M4").Value = 12.28762933333333
$ws.Range("N4").Value = 36.862888
$ws.Range("O4").Value = 0.4769116637062769
$ws.Range("P4").Value = 0.4769116637062769
$ws.Range("Q4").Value = 12.97838660805956
$ws.Range("R4").Value = 116.805479472536
$ws.Range("S4").Value = 0.4769116637062769
$ws.Range("T4").Value = 0.4769116637062769

# Row 5
$ws.Range("M5").Value = 2.085501
$ws.Range("N5").Value = 6.256503
$ws.Range("O5").Value = 0.08094317663644024
$ws.Range("P5").Value = 0.08094317663644023
$ws.Range("Q5").Value = 2.202738829049
$ws.Range("R5").Value = 19.824649461441
$ws.Range("S5").Value = 0.08094317663644024
$ws.Range("T5").Value = 0.08094317663644023
